# Commit: "added Account Lockout Test in Test package"
# Fix the typo "Account Lock Out Test" -> "Account Lockout Test"
# on the LoginCredentials sheet (cell D3), and move the active
# selection to D5 to match the saved workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginCredentials")

$ws.Range("D3").Value = "Account Lockout Test"

$ws.Range("D5").Select()
